$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 13 (shifts old rows 13-24 down to 14-25) ---
$ws.Rows.Item(13).Insert()

# The inserted row 13 picked up stray formatting in column A from the row above;
# clear that cell completely so the row only carries data in B/C.
$ws.Cells.Item(13, 1).Clear()

# --- Row 10: Objetivos: value (Portuguese objectives text) ---
$ws.Cells.Item(10, 2).Value = "Proporcionar ao aluno um conhecimento básico da mecânica dos corpos rígidos com ênfase na Estática dos Sólidos"
$ws.Cells.Item(10, 3).Value = "Proporcionar ao aluno um conhecimento básico da mecânica dos corpos rígidos com ênfase na Estática dos Sólidos"

# --- Row 13 (new): Docentes responsáveis value ---
$ws.Cells.Item(13, 2).Value = "5840650 - Janaína Ferreira Batista"
$ws.Cells.Item(13, 3).Value = "5840650 - Janaína Ferreira Batista"

# The newly inserted B13 cell inherited the wrong (bold) format because of the
# pre-existing column A/B width overlap; fix it up by copying the normal
# wrapped-text format from the equivalent cell below (row 14, same column).
$ws.Cells.Item(14, 2).Copy()
$ws.Cells.Item(13, 2).PasteSpecial(-4122)
$ws.Cells.Item(14, 3).Copy()
$ws.Cells.Item(13, 3).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 14: Programa resumido value ---
$ws.Cells.Item(14, 2).Value = "Estática das Partículas, Estática de Corpos Rígidos, Equilíbrio de Corpos Rígidos, Análise de Estruturas, Forças Distribuídas, Forças em Vigas."
$ws.Cells.Item(14, 3).Value = "Estática das Partículas, Estática de Corpos Rígidos, Equilíbrio de Corpos Rígidos, Análise de Estruturas, Forças Distribuídas, Forças em Vigas."

# --- Row 16: Programa value ---
$ws.Cells.Item(16, 2).Value = "Programa: 1) Estática de partículas: Vetores. Resultante de várias forças concorrentes. Equilíbrio de uma partícula. 2) Estática de Corpos Rígidos: Conceito de corpo rígido. Forças externas e forças internas. Forças equivalentes. Momento de uma força com relação a um ponto. Sistemas equivalentes de forças. Diagrama de corpo livre.3) Equilíbrio de corpos rígidos: Reações de apoios e conexões para uma estrutura 2D. Equilíbrio de um corpo rígido em 2D. Reações de apoios e conexões para uma estrutura 3D. Equilíbrio de um corpo rígido em 3D. 4) Análise de Estruturas: Treliças: Definições. Treliça simples. Análise de treliças pelo método dos nós. Análise de treliças pelo método das seções. Estruturas: estruturas que contêm elementos sujeitos a ação de múltiplas forças, transmissão e modificação de forças.5) Forças Distribuídas: Centróides e baricentros de áreas, linhas e volumes; teoremas de Guldinus-Pappus; cargas distribuídas sobre vigas.6) Forças em Vigas: Diagramas de cisalhamento e momento fletor para uma viga carregada."
$ws.Cells.Item(16, 3).Value = "Programa: 1) Estática de partículas: Vetores. Resultante de várias forças concorrentes. Equilíbrio de uma partícula. 2) Estática de Corpos Rígidos: Conceito de corpo rígido. Forças externas e forças internas. Forças equivalentes. Momento de uma força com relação a um ponto. Sistemas equivalentes de forças. Diagrama de corpo livre.3) Equilíbrio de corpos rígidos: Reações de apoios e conexões para uma estrutura 2D. Equilíbrio de um corpo rígido em 2D. Reações de apoios e conexões para uma estrutura 3D. Equilíbrio de um corpo rígido em 3D. 4) Análise de Estruturas: Treliças: Definições. Treliça simples. Análise de treliças pelo método dos nós. Análise de treliças pelo método das seções. Estruturas: estruturas que contêm elementos sujeitos a ação de múltiplas forças, transmissão e modificação de forças.5) Forças Distribuídas: Centróides e baricentros de áreas, linhas e volumes; teoremas de Guldinus-Pappus; cargas distribuídas sobre vigas.6) Forças em Vigas: Diagramas de cisalhamento e momento fletor para uma viga carregada."

# --- Row 19: Método value ---
$ws.Cells.Item(19, 2).Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Cells.Item(19, 3).Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# --- Row 20: Critério value ---
$ws.Cells.Item(20, 2).Value = "NF≥ 5,0."
$ws.Cells.Item(20, 3).Value = "NF≥ 5,0."

# --- Row 21: Norma de recuperação value ---
$ws.Cells.Item(21, 2).Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Cells.Item(21, 3).Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."

# --- Row 22: Bibliografia value ---
$ws.Cells.Item(22, 2).Value = "1. BEER, Ferdinand Pierre, ; JOHNSTON, E. Russel.; Eisenberg, Elliot R., Mecânica vetorial para engenheiros: Estática.  Mc Graw Hill (2011).2. HIBBELER, R.C. Mecânica para engenharia, Vol. 1: estática, Pearson Prentice Hall (2005).3. MERIAM J. L. ; KRAIGE, L. G., Mecânica, estática, Livros Técnicos e Científicos Editora (2004)."
$ws.Cells.Item(22, 3).Value = "1. BEER, Ferdinand Pierre, ; JOHNSTON, E. Russel.; Eisenberg, Elliot R., Mecânica vetorial para engenheiros: Estática.  Mc Graw Hill (2011).2. HIBBELER, R.C. Mecânica para engenharia, Vol. 1: estática, Pearson Prentice Hall (2005).3. MERIAM J. L. ; KRAIGE, L. G., Mecânica, estática, Livros Técnicos e Científicos Editora (2004)."

Write-Host "edit complete"
